# Generate Report for Handoff
#
# The localization-status report is regenerated: the source markdown file's
# GUID changed from 2b7f256f-28f6-45b7-b8f4-6bb54d494183 to
# 13068ee2-4fca-43b2-aadc-9abf8db79d87, the handoff bundle hash changed from
# b769ee4667bfcb04df0ff16b2d0ea89bdb1195b2 to
# 3f38102997147d12a3d76a1a6081d9bc16360811, and the handoff timestamps were
# refreshed.

$wb = $excel.ActiveWorkbook

$oldGuid = "2b7f256f-28f6-45b7-b8f4-6bb54d494183"
$newGuid = "13068ee2-4fca-43b2-aadc-9abf8db79d87"

# --- Overview sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = "$newGuid.md"

# --- zh-cn sheet ------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("C2").Value = "$newGuid.3f38102997147d12a3d76a1a6081d9bc16360811.zh-cn.xlf"
$ws.Range("D2").Value = "2016-03-10 12:07:19"

# --- de-de sheet ------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("C2").Value = "$newGuid.3f38102997147d12a3d76a1a6081d9bc16360811.de-de.xlf"
$ws.Range("D2").Value = "2016-03-10 12:07:23"
